$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'bjj compression shorts'
$ws.Cells.Item(2, 1).Value = 'compression running pants'
$ws.Cells.Item(3, 1).Value = 'sit pad'
$ws.Cells.Item(4, 1).Value = 'paintball pants padded'
$ws.Cells.Item(5, 1).Value = 'knee pads breathable'
$ws.Cells.Item(6, 1).Value = 'basketball compression tights for women'
$ws.Cells.Item(7, 1).Value = 'elastic band black mountain'
$ws.Cells.Item(8, 1).Value = 'lacrosse pads youth boys'
$ws.Cells.Item(9, 1).Value = 'position pad'
$ws.Cells.Item(10, 1).Value = 'knee length tights'
$ws.Cells.Item(11, 1).Value = 'capri pants men'
$ws.Cells.Item(12, 1).Value = 'knee pads volleyball black'
$ws.Cells.Item(13, 1).Value = 'softball sliding pants women'
$ws.Cells.Item(14, 1).Value = '6ft basketball'
$ws.Cells.Item(15, 1).Value = 'basketball shorts and pants'
$ws.Cells.Item(16, 1).Value = 'work pants for men construction knee pads'
$ws.Cells.Item(17, 1).Value = 'sliding shorts women softball'
$ws.Cells.Item(18, 1).Value = 'fight shorts wrestling'
$ws.Cells.Item(19, 1).Value = 'youth mesh leggings'
$ws.Cells.Item(20, 1).Value = 'above the knee basketball shorts'
$ws.Cells.Item(21, 1).Value = 'mens yoga leggings'
$ws.Cells.Item(22, 1).Value = 'weight lifting pants for men'
$ws.Cells.Item(23, 1).Value = 'cheap volleyball knee pads'
$ws.Cells.Item(24, 1).Value = 'compression spandex'
$ws.Cells.Item(25, 1).Value = 'yoga pants compression'
$ws.Cells.Item(26, 1).Value = 'mens above the knee shorts'
$ws.Cells.Item(27, 1).Value = 'mens running compression'
$ws.Cells.Item(28, 1).Value = 'cycling pants mens'
$ws.Cells.Item(29, 1).Value = 'knee sleeves basketball men'
$ws.Cells.Item(30, 1).Value = 'softball gear for girls'
$ws.Cells.Item(31, 1).Value = 'cold knee compression'
$ws.Cells.Item(32, 1).Value = 'youth padded compression shorts'
$ws.Cells.Item(33, 1).Value = 'yoga pants for men'
$ws.Cells.Item(34, 1).Value = 'mens spandex tights'
$ws.Cells.Item(35, 1).Value = 'softball protective gear'
$ws.Cells.Item(36, 1).Value = 'soccer sliding shorts'
$ws.Cells.Item(37, 1).Value = 'compression baseball shorts'
$ws.Cells.Item(38, 1).Value = 'long shorts for men below knee'
$ws.Cells.Item(39, 1).Value = 'padded leggings for cycling'
$ws.Cells.Item(40, 1).Value = 'padded volleyball shorts'
$ws.Cells.Item(41, 1).Value = 'hex squat'
$ws.Cells.Item(42, 1).Value = 'youth padded sliding shorts'
$ws.Cells.Item(43, 1).Value = 'knee sleeves bjj'
$ws.Cells.Item(44, 1).Value = 'football pants pads adult'
$ws.Cells.Item(45, 1).Value = 'work pants knee'
$ws.Cells.Item(46, 1).Value = 'cold compression knee'
$ws.Cells.Item(47, 1).Value = '5 pad football girdle'
$ws.Cells.Item(48, 1).Value = 'wrestling sleeve youth'
$ws.Cells.Item(49, 1).Value = 'compression sports pants'
$ws.Cells.Item(50, 1).Value = 'basketball tights for girls'
$ws.Cells.Item(51, 1).Value = 'water pants'
$ws.Cells.Item(52, 1).Value = 'spandex tights men'
$ws.Cells.Item(53, 1).Value = 'boys compression pants black'
$ws.Cells.Item(54, 1).Value = 'hockey tights'
$ws.Cells.Item(55, 1).Value = 'youth hockey compression pants'
$ws.Cells.Item(56, 1).Value = 'men leggings compression'
$ws.Cells.Item(57, 1).Value = 'wrestling kneepads'
$ws.Cells.Item(58, 1).Value = 'kneeling pad gym'
$ws.Cells.Item(59, 1).Value = 'guard shorts'
$ws.Cells.Item(60, 1).Value = 'padded compression shorts men'
$ws.Cells.Item(61, 1).Value = 'softball pants youth'
$ws.Cells.Item(62, 1).Value = 'spandex basketball shorts'
$ws.Cells.Item(63, 1).Value = 'compression shorts men 5 pack'
$ws.Cells.Item(64, 1).Value = 'shorts for men below knee'
$ws.Cells.Item(65, 1).Value = 'mens gym leggings'
$ws.Cells.Item(66, 1).Value = 'compression running leggings'
$ws.Cells.Item(67, 1).Value = 'black mens basketball shorts'
$ws.Cells.Item(68, 1).Value = 'knee pads impact'
$ws.Cells.Item(69, 1).Value = 'paintball pads'
$ws.Cells.Item(70, 1).Value = 'boys compression'
$ws.Cells.Item(71, 1).Value = 'mens volleyball kneepads'
$ws.Cells.Item(72, 1).Value = 'yoga knee pads'
$ws.Cells.Item(73, 1).Value = 'knee work pads'
$ws.Cells.Item(74, 1).Value = 'running capri'
$ws.Cells.Item(75, 1).Value = 'paintball pants for men'
$ws.Cells.Item(76, 1).Value = 'kneepad youth'
$ws.Cells.Item(77, 1).Value = 'polyester capri pants'
$ws.Cells.Item(78, 1).Value = 'man capri pants'
$ws.Cells.Item(79, 1).Value = 'indoor baseball'
$ws.Cells.Item(80, 1).Value = 'softball compression sleeve'
$ws.Cells.Item(81, 1).Value = 'male pads'
$ws.Cells.Item(82, 1).Value = 'high five girls softball pants'
$ws.Cells.Item(83, 1).Value = 'outdoor hockey pants'
$ws.Cells.Item(84, 1).Value = 'basketball floor'
$ws.Cells.Item(85, 1).Value = 'basketball knee sleeve black'
$ws.Cells.Item(86, 1).Value = 'long shorts for men below knee sports'
$ws.Cells.Item(87, 1).Value = 'knee pads for adults'
$ws.Cells.Item(88, 1).Value = 'hockey leggings'
$ws.Cells.Item(89, 1).Value = 'volleyball long knee pads'
$ws.Cells.Item(90, 1).Value = 'lacrosse shorts mens'
$ws.Cells.Item(91, 1).Value = 'mens tights with pouch'
$ws.Cells.Item(92, 1).Value = 'black short baseball pants'
$ws.Cells.Item(93, 1).Value = 'lightweight knee pads'
$ws.Cells.Item(94, 1).Value = 'mens compression pants cold'
$ws.Cells.Item(95, 1).Value = 'knee shorts'
$ws.Cells.Item(96, 1).Value = 'girls sliding pants'
$ws.Cells.Item(97, 1).Value = 'knee pads for work for men'
$ws.Cells.Item(98, 1).Value = 'youth padded leg sleeves for basketball'
$ws.Cells.Item(99, 1).Value = 'gym knee compression'
$ws.Cells.Item(100, 1).Value = 'compression football girdle'
